$d = $word.ActiveDocument

# 1) Highlight (green) the first run of the "Un compratore può acquistare..." paragraph
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Un compratore può acquistare energia da uno o più venditori per un giorno. Per una fascia oraria un compratore può comprare solo da un produttore.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Font.HighlightColorIndex = 4
}

# 2) Strike-through the "Dare la possibilità ad un consumatore di riservare..." list item
#    (both the run text and the paragraph mark get the strike formatting)
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Dare la possibilità ad un consumatore di riservare uno slot per il giorno seguente in una fascia oraria. L’acquisto minimo è di 0.1kWh", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $para2.Range.Font.StrikeThrough = 1
}

# 3) Highlight (green) the "Dare la possibilità ad un consumatore di modificare..." list item
#    (covers all runs plus the paragraph mark, preserving existing bold/underline formatting)
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Dare la possibilità ad un consumatore di modificare (anche ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $para3 = $rng3.Paragraphs(1)
    $para3.Range.Font.HighlightColorIndex = 4
}
